$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.783.79"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "3.804.96"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "597.07"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "167.50"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "0.161"
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "6.29"
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "35.98"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").Value = "4.436.49"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "3.808.39"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("D17").Value = "67.843.42"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "7.08"
$ws.Range("E18").Value = "  +1.25%  "
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").Value = "461.56"
$ws.Range("E20").Value = "  +0.46%  "
$ws.Range("D21").Value = "9.90"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("D22").Value = "0.699"
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "83.35"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").Value = "12.06"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("E26").Value = "  -1.87%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").Value = "3.951.01"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "7.36"
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").Value = "29.49"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D34").Value = "9.04"
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Value = "3.745.29"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.100"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.40"
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.138"
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "5.77"
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "48.06"
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").Value = "0.301"
$ws.Range("E45").Value = "  +0.91%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").Value = "42.75"
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "1.38"
$ws.Range("E47").Value = "  +8.91%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").Value = "8.33"
$ws.Range("E48").Value = "  -0.49%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "27.37"
$ws.Range("E49").Value = "  +8.53%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "147.84"
$ws.Range("E50").Value = "  +0.17%  "
$ws.Range("B51").Value = "Bittensor"
$ws.Range("C51").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D51").Value = "393.66"
$ws.Range("E51").Value = "  +0.20%  "
